$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("B6").Value = "'2024-05-24"
$ws.Range("C6").Value = "北京·2024国际收藏玩具与艺术创意展览会"
$ws.Range("D6").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E6").Value = "2024.05.24 09:30-05.26 18:00"
$ws.Range("F6").Value = 65
$ws.Range("G6").Value = 120
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84698"
$ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/W205PYEt1713855610417.jpeg"
$ws.Range("B7").Value = "'2024-05-25"
$ws.Range("C7").Value = "北京·LookLook剧情式沉浸游戏互动动漫嘉年华"
$ws.Range("D7").Value = "东村文化创意产业园A1-2 五道杠实景片场"
$ws.Range("E7").Value = "2024.05.25 09:30-05.26 17:30"
$ws.Range("F7").Value = 870
$ws.Range("G7").Value = 72
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=84741"
$ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202404/ytcuhFxO1713859439829.jpeg"
$ws.Range("F14").Value = 1570
$ws.Range("F15").Value = 7227
$ws.Range("F17").Value = 7379
$ws.Range("F19").Value = 21
$ws.Range("F20").Value = 5750
$ws.Range("F21").Value = 5750
$ws.Range("F22").Value = 3156
$ws.Range("F23").Value = 3535
$ws.Range("F26").Value = 253
$ws.Range("F28").Value = 1980
$ws.Range("F30").Value = 321
$ws.Range("F31").Value = 893
$ws.Range("F32").Value = 240
$ws.Range("F33").Value = 511
$ws.Range("F34").Value = 49
$ws.Range("F35").Value = 2498
$ws.Range("F36").Value = 1311
$ws.Range("F37").Value = 2948
$ws.Range("F38").Value = 91
$ws.Range("F39").Value = 25
$ws.Range("F40").Value = 182
$ws.Range("F41").Value = 430
$ws.Range("F42").Value = 1147
$ws.Range("F44").Value = 501
# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = "不可售"
$ws.Range("F10").Value = 36
$ws.Range("F11").Value = 376
$ws.Range("F15").Value = 4
$ws.Range("F17").Value = 57
# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 94
# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("B7").Value = "'2024-05-20"
$ws.Range("C7").Value = "北京·奇妙嘿夜烛光音乐会“不能说的秘密”海洋主题"
$ws.Range("D7").Value = "万寿路街道复兴路69号五棵松万达广场6层 格乐丽雅（中国）艺术中心"
$ws.Range("E7").Value = "2024.05.20 18:30-05.20 21:50"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 438
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=84379"
$ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202404/lqa6YTqQ1713252023331.jpeg"
$ws.Range("B8").Value = "'2024-05-24"
$ws.Range("C8").Value = "北京·2024国际收藏玩具与艺术创意展览会"
$ws.Range("D8").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E8").Value = "2024.05.24 09:30-05.26 18:00"
$ws.Range("F8").Value = 65
$ws.Range("G8").Value = 120
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84698"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202404/W205PYEt1713855610417.jpeg"
$ws.Range("B9").Value = "'2024-05-25"
$ws.Range("C9").Value = "北京·LookLook剧情式沉浸游戏互动动漫嘉年华"
$ws.Range("D9").Value = "东村文化创意产业园A1-2 五道杠实景片场"
$ws.Range("E9").Value = "2024.05.25 09:30-05.26 17:30"
$ws.Range("F9").Value = 870
$ws.Range("G9").Value = 72
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84741"
$ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/ytcuhFxO1713859439829.jpeg"
$ws.Range("F14").Value = 94
$ws.Range("F16").Value = 1570
$ws.Range("F19").Value = 36
$ws.Range("F20").Value = 7226
$ws.Range("F22").Value = 7379
$ws.Range("F23").Value = 21
$ws.Range("F24").Value = 5750
$ws.Range("F25").Value = 5750
$ws.Range("F26").Value = 3156
$ws.Range("F27").Value = 3535
$ws.Range("F30").Value = 253
$ws.Range("F32").Value = 1980
$ws.Range("F34").Value = 57
$ws.Range("F35").Value = 321
$ws.Range("F36").Value = 893
$ws.Range("F37").Value = 511
$ws.Range("F38").Value = 49
$ws.Range("F39").Value = 2498
$ws.Range("F40").Value = 1311
$ws.Range("F42").Value = 2948
$ws.Range("F43").Value = 91
$ws.Range("F44").Value = 25
$ws.Range("F45").Value = 182
$ws.Range("F47").Value = 430
$ws.Range("F48").Value = 1147
$ws.Range("F50").Value = 501
